$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 415.16666
$ws.Range("J17").Value = 415.16666
$ws.Range("L17").Value = 1245.49998
$ws.Range("N17").Value = -1581.49998

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H33").Value = 1818.8889
$ws.Range("I33").Value = 2313.2856
$ws.Range("K33").Value = 2313.2856
$ws.Range("M33").Value = -2084.2856

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H69").Value = 9661.143
$ws.Range("I69").Value = 6537.6665
$ws.Range("K69").Value = 19612.9995
$ws.Range("M69").Value = -18738.9995

$ws.Range("H72").Value = 9661.143
$ws.Range("I72").Value = 6537.6665
$ws.Range("K72").Value = 58838.9985
$ws.Range("M72").Value = -54470.9985

$ws.Range("H98").Value = 3229.8
$ws.Range("I98").Value = 2849.182
$ws.Range("J98").Value = 3695
$ws.Range("K98").Value = 2849.182
$ws.Range("L98").Value = 3695
$ws.Range("M98").Value = -1351.182
$ws.Range("N98").Value = -6691

$ws.Range("H100").Value = 4172.3335
$ws.Range("I100").Value = 2509.5557
$ws.Range("J100").Value = 6666.5
$ws.Range("K100").Value = 2509.5557
$ws.Range("L100").Value = 6666.5
$ws.Range("M100").Value = -1968.5557
$ws.Range("N100").Value = -7748.5

$ws.Range("H122").Value = 3229.8
$ws.Range("I122").Value = 2849.182
$ws.Range("J122").Value = 3695
$ws.Range("K122").Value = 8547.545999999998
$ws.Range("L122").Value = 11085
$ws.Range("M122").Value = -6097.545999999998
$ws.Range("N122").Value = -15985

$ws.Range("H135").Value = 1443.5
$ws.Range("I135").Value = 1512.8148
$ws.Range("J135").Value = 1235.5555
$ws.Range("K135").Value = 13615.3332
$ws.Range("L135").Value = 11119.9995
$ws.Range("M135").Value = -11080.3332
$ws.Range("N135").Value = -16189.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1295.2322
$ws.Range("I32").Value = 1311.9623
$ws.Range("K32").Value = 1311.9623
$ws.Range("M32").Value = -1024.9623

$ws.Range("H61").Value = 2293.862
$ws.Range("I61").Value = 2366.423
$ws.Range("J61").Value = 1665
$ws.Range("K61").Value = 2366.423
$ws.Range("L61").Value = 1665
$ws.Range("M61").Value = -2154.423
$ws.Range("N61").Value = -2089

$ws.Range("H74").Value = 426300.25
$ws.Range("I74").Value = 455962.62
$ws.Range("K74").Value = 455962.62
$ws.Range("M74").Value = -455088.62

$ws.Range("H77").Value = 426300.25
$ws.Range("I77").Value = 455962.62
$ws.Range("K77").Value = 2279813.1
$ws.Range("M77").Value = -2275445.1

$ws.Range("H114").Value = 90398
$ws.Range("J114").Value = 90398
$ws.Range("L114").Value = 90398
$ws.Range("N114").Value = -99076

$ws.Range("H122").Value = 27782322
$ws.Range("I122").Value = 50003144
$ws.Range("J122").Value = 6293.0625
$ws.Range("K122").Value = 150009432
$ws.Range("L122").Value = 18879.1875
$ws.Range("M122").Value = -150006982
$ws.Range("N122").Value = -23779.1875

$ws.Range("H132").Value = 194046.6
$ws.Range("I132").Value = 219838.2
$ws.Range("K132").Value = 659514.6000000001
$ws.Range("M132").Value = -656984.6000000001

$ws.Range("H136").Value = 2293.862
$ws.Range("I136").Value = 2366.423
$ws.Range("J136").Value = 1665
$ws.Range("K136").Value = 7099.268999999999
$ws.Range("L136").Value = 4995
$ws.Range("M136").Value = -4549.268999999999
$ws.Range("N136").Value = -10095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2319.9285
$ws.Range("I20").Value = 1335.625
$ws.Range("K20").Value = 1335.625
$ws.Range("M20").Value = -1088.625

$ws.Range("H86").Value = 812190.5
$ws.Range("I86").Value = 1135586.4
$ws.Range("J86").Value = 3700.8333
$ws.Range("K86").Value = 1135586.4
$ws.Range("L86").Value = 3700.8333
$ws.Range("M86").Value = -1134463.4
$ws.Range("N86").Value = -5946.8333

$ws.Range("H89").Value = 812190.5
$ws.Range("I89").Value = 1135586.4
$ws.Range("J89").Value = 3700.8333
$ws.Range("K89").Value = 5677932
$ws.Range("L89").Value = 18504.1665
$ws.Range("M89").Value = -5672316
$ws.Range("N89").Value = -29736.1665

$ws.Range("H107").Value = 1879.909
$ws.Range("I107").Value = 1742.1111
$ws.Range("K107").Value = 1742.1111
$ws.Range("M107").Value = 177.8888999999999

$ws.Range("H134").Value = 32148.865
$ws.Range("I134").Value = 1470.2222
$ws.Range("K134").Value = 4410.6666
$ws.Range("M134").Value = -1875.6666

$ws.Range("H135").Value = 60193.75
$ws.Range("J135").Value = 60193.75
$ws.Range("L135").Value = 60193.75
$ws.Range("N135").Value = -70333.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7846.75
$ws.Range("I58").Value = 2526.5715
$ws.Range("K58").Value = 2526.5715
$ws.Range("M58").Value = -2323.5715

$ws.Range("H86").Value = 11084.857
$ws.Range("I86").Value = 3984.4285
$ws.Range("J86").Value = 18185.285
$ws.Range("K86").Value = 3984.4285
$ws.Range("L86").Value = 18185.285
$ws.Range("M86").Value = -2861.4285
$ws.Range("N86").Value = -20431.285

$ws.Range("H89").Value = 11084.857
$ws.Range("I89").Value = 3984.4285
$ws.Range("J89").Value = 18185.285
$ws.Range("K89").Value = 19922.1425
$ws.Range("L89").Value = 90926.425
$ws.Range("M89").Value = -14306.1425
$ws.Range("N89").Value = -102158.425

$ws.Range("H94").Value = 648.75
$ws.Range("J94").Value = 745.1
$ws.Range("L94").Value = 745.1
$ws.Range("N94").Value = -1647.1

$ws.Range("H107").Value = 1032.8846
$ws.Range("I107").Value = 692.05554
$ws.Range("J107").Value = 1799.75
$ws.Range("K107").Value = 692.05554
$ws.Range("L107").Value = 1799.75
$ws.Range("M107").Value = 1227.94446
$ws.Range("N107").Value = -5639.75

$ws.Range("H132").Value = 4193.769
$ws.Range("I132").Value = 3120.3333
$ws.Range("J132").Value = 5113.857
$ws.Range("K132").Value = 9360.999899999999
$ws.Range("L132").Value = 15341.571
$ws.Range("M132").Value = -6830.999899999999
$ws.Range("N132").Value = -20401.571

$ws.Range("H136").Value = 7846.75
$ws.Range("I136").Value = 2526.5715
$ws.Range("K136").Value = 7579.7145
$ws.Range("M136").Value = -5029.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 15117592
$ws.Range("I4").Value = 895902.6
$ws.Range("J4").Value = 120358100
$ws.Range("K4").Value = 2687707.8
$ws.Range("L4").Value = 361074300
$ws.Range("M4").Value = -2687595.8
$ws.Range("N4").Value = -361074524

$ws.Range("H34").Value = 8687.8125
$ws.Range("I34").Value = 158.42857
$ws.Range("J34").Value = 15321.777
$ws.Range("K34").Value = 475.28571
$ws.Range("L34").Value = 45965.331
$ws.Range("M34").Value = -391.28571
$ws.Range("N34").Value = -46133.331

$ws.Range("H122").Value = 10451086
$ws.Range("J122").Value = 1233.2858
$ws.Range("L122").Value = 11099.5722
$ws.Range("N122").Value = -15999.5722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 557.3333
$ws.Range("J107").Value = 813.6
$ws.Range("L107").Value = 813.6
$ws.Range("N107").Value = -4653.6

$ws.Range("H122").Value = 427458.53
$ws.Range("I122").Value = 554373.75
$ws.Range("K122").Value = 1663121.25
$ws.Range("M122").Value = -1660671.25

$ws.Range("H132").Value = 66836.62
$ws.Range("I132").Value = 25127.133
$ws.Range("J132").Value = 171110.33
$ws.Range("K132").Value = 75381.399
$ws.Range("L132").Value = 513330.99
$ws.Range("M132").Value = -72851.399
$ws.Range("N132").Value = -518390.99

$ws.Range("H134").Value = 87161
$ws.Range("J134").Value = 87161
$ws.Range("L134").Value = 261483
$ws.Range("N134").Value = -266553

$ws.Range("H136").Value = 74550.22
$ws.Range("J136").Value = 74550.22
$ws.Range("L136").Value = 223650.66
$ws.Range("N136").Value = -228750.66

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H46").Value = 2123
$ws.Range("I46").Value = 1849.5555
$ws.Range("J46").Value = 2511.5789
$ws.Range("K46").Value = 1849.5555
$ws.Range("L46").Value = 2511.5789
$ws.Range("M46").Value = -1661.5555
$ws.Range("N46").Value = -2887.5789

$ws.Range("H100").Value = 71786.875
$ws.Range("I100").Value = 77220.71000000001
$ws.Range("K100").Value = 77220.71000000001
$ws.Range("M100").Value = -76679.71000000001

$ws.Range("H132").Value = 5252.0605
$ws.Range("I132").Value = 4437.696
$ws.Range("J132").Value = 7125.1
$ws.Range("K132").Value = 13313.088
$ws.Range("L132").Value = 21375.3
$ws.Range("M132").Value = -10783.088
$ws.Range("N132").Value = -26435.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1286.5
$ws.Range("I107").Value = 1718.75
$ws.Range("J107").Value = 998.3333
$ws.Range("K107").Value = 5156.25
$ws.Range("L107").Value = 2994.9999
$ws.Range("M107").Value = -3236.25
$ws.Range("N107").Value = -6834.9999

$ws.Range("H113").Value = 1466.75
$ws.Range("I113").Value = 1378.579
$ws.Range("K113").Value = 4135.737
$ws.Range("M113").Value = -1965.737

$ws.Range("H136").Value = 45401.1
$ws.Range("I136").Value = 2145.0256
$ws.Range("K136").Value = 6435.0768
$ws.Range("M136").Value = -3885.0768

$ws.Range("H138").Value = 73660
$ws.Range("J138").Value = 73660
$ws.Range("L138").Value = 73660
$ws.Range("N138").Value = -83940
